$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Bill of Materials" part numbers/models in column B
# (order matters for shared-string table layout, matching the original author's edit order)
$ws.Range("B4").Value = "HC-SR501"
$ws.Range("B2").Value = "LDR"
$ws.Range("B7").Value = "NA"
$ws.Range("B6").Value = "ACS712"
$ws.Range("B5").Value = "SRD-05VDC-SL-C"

# Update the active selection to match the author's final cursor position
$ws.Range("B5").Select()
